# Update the "marksheet" summary values (Corr/total marks) on the quiz sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: Right count 3 -> 5
$ws.Range("B11").Value = 5

# Total row: Right marks total 66 -> 110
$ws.Range("B12").Value = 110

# Total row: Correct/Total marks text "63/84" -> "110/140"
$ws.Range("E12").Value = "110/140"
